# Correção documentos diversos (Feedback do Professor)
# Update the work-items descriptions on the "Lista de Itens de Trabalho" sheet
# to be more descriptive (mentioning "Sistema de Rastreamento"), and move the
# selection to A14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Itens de Trabalho")
$ws.Activate()

$ws.Range("A3").Value = "Realizar Análise de Sistema"
$ws.Range("A2").Value = "Requisitos do Sistema de Rastreamento"
$ws.Range("A4").Value = "Criar Projeto do Sistema de Rastreamento"
$ws.Range("A5").Value = "Implementação do Sistema de Rastreamento"
$ws.Range("A6").Value = "Realizar Testes no Sistema de Rastreamento"
$ws.Range("A7").Value = "Implantação do Sistema de Rastreamento"
$ws.Range("A8").Value = "Criar tela de Login do Sistema de Rastreamento"
$ws.Range("A9").Value = "Criar Tela de  Localização de Veículos"

$ws.Range("A14").Select()
